$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value is a plain number-looking string (single
# decimal point, e.g. "15.00", "0.0800") must be forced to Text formatting
# before assignment, otherwise Excel auto-converts them to numbers and
# mangles trailing zeros / switches to scientific notation. The format is
# reset back to Normal immediately after the write so no stray formatting
# is left behind. Values already containing two dots (e.g. "37.017.93")
# are never parsed as numbers by Excel, so they do not need this treatment.

$ws.Range("D2").Value = '37.017.93'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '2.058.44'
$ws.Range("E3").Value = '  -2.10%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.669'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.67%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.28'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("E10").Value = '  +0.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0799'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.89%  '

$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.04%  '

$ws.Range("D14").Value = '2.361.57'
$ws.Range("E14").Value = '  -2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.810'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.25%  '

$ws.Range("D17").Value = '2.064.52'
$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("D18").Value = '36.927.77'
$ws.Range("E18").Value = '  -0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.03'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = '0.0₃0925'
$ws.Range("E20").Value = '  +11.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("E25").Value = '  -4.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.56%  '

$ws.Range("E29").Value = '  -1.59%  '

$ws.Range("E30").Value = '  +0.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0627'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.39%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.21%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0873'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.61%  '

$ws.Range("E37").Value = '  -6.12%  '

$ws.Range("E38").Value = '  -1.55%  '

$ws.Range("E39").Value = '  -0.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +21.63%  '

$ws.Range("B41").Value = 'FTXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +63.08%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.44%  '

$ws.Range("E43").Value = '  -0.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '

$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.74%  '

$ws.Range("D49").Value = '1.297.36'
$ws.Range("E49").Value = '  -4.44%  '

$ws.Range("E50").Value = '  -0.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.24%  '
